# RC20XX File Transfer instructions.docx - apply commit "Update File transfer docs"
#
# Strategy: Word's proofing engine normally splits text into multiple
# <w:r> runs around <w:proofErr> marks (gramStart/gramEnd, spellStart/
# spellEnd) as a side effect of spell/grammar checking while typing -
# that isn't reachable from the exposed Word object model here, so we
# rebuild the affected paragraphs directly via Range.InsertXML, which
# accepts a raw OOXML fragment and inserts it at/after the given Range.
# We locate each target paragraph with Find, delete its Range, then
# insert the replacement run/proofErr structure in its place.

$d = $word.ActiveDocument

function Replace-ParaContent {
    param(
        [string]$FindText,
        [string]$InnerXml
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $FindText"
    }
    $rng.Delete() | Out-Null
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $InnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml) | Out-Null
}

# 1. "Start? Is displayed in Terminal " -> split around "Terminal" (gramStart/gramEnd)
$inner = '<w:p><w:r><w:t xml:space="preserve">Start? Is displayed in </w:t></w:r>' +
         '<w:proofErr w:type="gramStart"/><w:r><w:t>Terminal</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
Replace-ParaContent "Start? Is displayed in Terminal " $inner

# 2. "On your PC, Close your terminal editor and Go to the directory With the Python Files"
$inner = '<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">On your PC, </w:t></w:r>' +
         '<w:proofErr w:type="gramStart"/><w:r><w:t>Close</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
         '<w:r><w:t xml:space="preserve"> your terminal editor and Go to the directory With the Python Files</w:t></w:r></w:p>'
Replace-ParaContent "On your PC, Close your terminal editor and Go to the directory With the Python Files" $inner

# 3. "Type LS-RC20xx.py COMPORT DriveLetter" -> keep "Type LS-RC20xx" + "." runs, split last run
$inner = '<w:p><w:r><w:t>Type LS-RC20xx</w:t></w:r><w:r><w:t>.</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">py COMPORT </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>DriveLetter</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Replace-ParaContent "Type LS-RC20xx.py COMPORT DriveLetter" $inner

# 4. "LS-RC20xx.py COM4 D" -> merge the six runs into one
$inner = '<w:p><w:r><w:t>LS-RC20xx.py COM4 D</w:t></w:r></w:p>'
Replace-ParaContent "LS-RC20xx.py COM4 D" $inner

# 5. "CopyFrom-RC20xx.py ComPort DriveLetter FileName"
$inner = '<w:p><w:r><w:t xml:space="preserve">CopyFrom-RC20xx.py </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>ComPort</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>DriveLetter</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>FileName</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Replace-ParaContent "CopyFrom-RC20xx.py ComPort DriveLetter FileName" $inner

# 6. "CopyFrom-RC20xx.py" + " COM4 D TESTMAND.BAS" -> merge into one run
$inner = '<w:p><w:r><w:t>CopyFrom-RC20xx.py COM4 D TESTMAND.BAS</w:t></w:r></w:p>'
Replace-ParaContent "CopyFrom-RC20xx.py COM4 D TESTMAND.BAS" $inner

# 7. "Copy" + "To" + "-RC20xx.py ComPort DriveLetter FileName" -> "CopyTo-RC20xx.py ..."
$inner = '<w:p><w:r><w:t xml:space="preserve">CopyTo-RC20xx.py </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>ComPort</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>DriveLetter</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>FileName</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Replace-ParaContent "CopyTo-RC20xx.py ComPort DriveLetter FileName" $inner

# 8. "RM" + "-RC20xx.py ComPort DriveLetter FileName"
$inner = '<w:p><w:r><w:t xml:space="preserve">RM-RC20xx.py </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>ComPort</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>DriveLetter</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>FileName</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Replace-ParaContent "RM-RC20xx.py ComPort DriveLetter FileName" $inner

# 9. "RM" + "-RC20xx.py" (+ unchanged " ", "COM4 D ", "TEST.TXT") -> merge first two runs
$inner = '<w:p><w:r><w:t>RM-RC20xx.py</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">COM4 D </w:t></w:r><w:r><w:t>TEST.TXT</w:t></w:r></w:p>'
Replace-ParaContent "RM-RC20xx.py COM4 D TEST.TXT" $inner

# 10. "To " + "cat (type)" + " a file from the RC20xx" -> merge into one run
$inner = '<w:p><w:r><w:t>To cat (type) a file from the RC20xx</w:t></w:r></w:p>'
Replace-ParaContent "To cat (type) a file from the RC20xx" $inner

# 11. "CAT" + "-RC20xx.py ComPort DriveLetter FileName"
$inner = '<w:p><w:r><w:t xml:space="preserve">CAT-RC20xx.py </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>ComPort</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>DriveLetter</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>FileName</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Replace-ParaContent "CAT-RC20xx.py ComPort DriveLetter FileName" $inner

# 12. "CAT" + "-RC20xx.py" (+ unchanged " ", "COM4 D ", "TESTMAND.BAS") -> merge first two runs
$inner = '<w:p><w:r><w:t>CAT-RC20xx.py</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">COM4 D </w:t></w:r><w:r><w:t>TESTMAND.BAS</w:t></w:r></w:p>'
Replace-ParaContent "CAT-RC20xx.py COM4 D TESTMAND.BAS" $inner

# 13. "EXIT-RC2040.py" + " COmPort" -> keep first run, split the rest
$inner = '<w:p><w:r><w:t>EXIT-RC2040.py</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:r><w:t>COmPort</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Replace-ParaContent "EXIT-RC2040.py COmPort" $inner

# 14. "RC2040 LED off" -> split around "off" (gramStart/gramEnd)
$inner = '<w:p><w:r><w:t xml:space="preserve">RC2040 LED </w:t></w:r>' +
         '<w:proofErr w:type="gramStart"/><w:r><w:t>off</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
Replace-ParaContent "RC2040 LED off" $inner

# 15. "And if you reconnect your terminal you are back to the console"
$inner = '<w:p><w:r><w:t xml:space="preserve">And if you reconnect your </w:t></w:r>' +
         '<w:proofErr w:type="gramStart"/><w:r><w:t>terminal</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
         '<w:r><w:t xml:space="preserve"> you are back to the console</w:t></w:r></w:p>'
Replace-ParaContent "And if you reconnect your terminal you are back to the console" $inner

# 16. Append new paragraphs at the end of the document (after the last
#     picture, before the trailing empty paragraph / sectPr).
$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
          '<w:p/>' +
          '<w:p/>' +
          '<w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
          '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">New </w:t></w:r></w:p>' +
          '<w:p><w:r><w:t xml:space="preserve">Program the RC2040 with a binary file. </w:t></w:r></w:p>' +
          '<w:p><w:r><w:t xml:space="preserve">After clicking the button. </w:t></w:r></w:p>' +
          '<w:p><w:r><w:t xml:space="preserve">Run the Python </w:t></w:r>' +
          '<w:proofErr w:type="gramStart"/><w:r><w:t>script</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
          '<w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' +
          '<w:p><w:r><w:t>Program-RC2040.py COM</w:t></w:r><w:r><w:t>4</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> 0000 24.bin</w:t></w:r></w:p>' +
          '<w:p><w:r><w:t xml:space="preserve">Where 0000 is the base Address and </w:t></w:r>' +
          '<w:proofErr w:type="gramStart"/><w:r><w:t>24.Bin</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
          '<w:r><w:t xml:space="preserve"> is  a raw binary file. </w:t></w:r></w:p>' +
          '<w:p><w:r><w:t xml:space="preserve">Exit the FFS mode and your RC2040 has the new code at the address you specified. </w:t></w:r></w:p>' +
          '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$target = $lastPara.Range
$target.Collapse(1) | Out-Null
$target.InsertXML($newXml) | Out-Null
